$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new column before column H (8) for pct_total_profit
$ws.Columns.Item(8).Insert()
# Insert new column before column J (10) for pct_total_volume (after total_volume, which is now column I)
$ws.Columns.Item(10).Insert()

# Set the new column widths to match the target layout
# (subtract 5/6 to compensate for the COM ColumnWidth -> stored-width offset)
$ws.Columns.Item(8).ColumnWidth = 20.166666666666668
$ws.Columns.Item(9).ColumnWidth = 19.166666666666668
$ws.Columns.Item(10).ColumnWidth = 19.166666666666668
$ws.Columns.Item(11).ColumnWidth = 45.166666666666664

# Header row 2 (EURUSD block)
$ws.Range("H2").Value = "pct_total_profit"
$ws.Range("J2").Value = "pct_total_volume"

# Header row 9 (GBPUSD block)
$ws.Range("H9").Value = "pct_total_profit"
$ws.Range("J9").Value = "pct_total_volume"

# Header row 16 (USDJPY block)
$ws.Range("H16").Value = "pct_total_profit"
$ws.Range("J16").Value = "pct_total_volume"

# Header row 23 (XAUUSD block)
$ws.Range("H23").Value = "pct_total_profit"
$ws.Range("J23").Value = "pct_total_volume"

# Data values - EURUSD block
$ws.Range("H4").Value = 40.38088209480964
$ws.Range("J4").Value = 39.9050663410838
$ws.Range("H5").Value = 58.83966121977033
$ws.Range("J5").Value = 22.11342858164777
$ws.Range("H6").Value = -8.011031082333684
$ws.Range("J6").Value = 24.41730584623592
$ws.Range("H7").Value = 8.790487767753721
$ws.Range("J7").Value = 13.5641992310325

# Data values - GBPUSD block
$ws.Range("H11").Value = 15.46061711908153
$ws.Range("J11").Value = 41.13726492411173
$ws.Range("H12").Value = 59.25583569534241
$ws.Range("J12").Value = 14.99113688190518
$ws.Range("H13").Value = 3.642030385105345
$ws.Range("J13").Value = 38.74859575221171
$ws.Range("H14").Value = 21.64151680047071
$ws.Range("J14").Value = 5.123002441771371

# Data values - USDJPY block
$ws.Range("H18").Value = 32.18919904839271
$ws.Range("J18").Value = 40.75402197017907
$ws.Range("H19").Value = 53.28339273949857
$ws.Range("J19").Value = 28.73028686822357
$ws.Range("H20").Value = 3.069190698721309
$ws.Range("J20").Value = 18.92621581632239
$ws.Range("H21").Value = 11.45821751338742
$ws.Range("J21").Value = 11.58947534527499

# Data values - XAUUSD block
$ws.Range("H25").Value = 52.10661977616129
$ws.Range("J25").Value = 40.4110727068207
$ws.Range("H26").Value = 62.42484337622932
$ws.Range("J26").Value = 25.03852608783536
$ws.Range("H27").Value = -13.86767213287978
$ws.Range("J27").Value = 23.50951214029568
$ws.Range("H28").Value = -0.6637910195108355
$ws.Range("J28").Value = 11.04088906504825
